$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from after the second "Surname" entry
#    (under "Supplier details:") to right after the "Wish list:" run.
#    Word only keeps one bookmark per name, so adding a new "_GoBack"
#    bookmark automatically removes the old one.
#
#    This host's Bookmarks.Add() mishandles a zero-length Range sitting
#    exactly on a paragraph-mark boundary (it silently falls back to
#    Paragraph(1)'s range), so instead we insert a temporary one
#    character placeholder right after the text, wrap a (non-collapsed)
#    bookmark around it, then delete the placeholder -- which correctly
#    collapses the bookmark in place.
# ---------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("Wish list:", $false, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $placeholderStart = $findRange.End
    $insertPoint = $d.Range($placeholderStart, $placeholderStart)
    $insertPoint.InsertAfter("X")

    $placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
    $d.Bookmarks.Add("_GoBack", $placeholderRange)

    $placeholderRange2 = $d.Range($placeholderStart, $placeholderStart + 1)
    $placeholderRange2.Delete()
}

# ---------------------------------------------------------------------
# 2) Colour the "Orders:" heading (and its paragraph mark) red, to
#    match "Wish list:".
# ---------------------------------------------------------------------
$ordersFindRange = $d.Content
$ordersFound = $ordersFindRange.Find.Execute("Orders:", $false, $false, $false, $false, $false, `
                                              $true, 1, $false, "", 0)
if ($ordersFound) {
    $ordersParagraph = $ordersFindRange.Paragraphs(1)
    $ordersParagraph.Range.Font.Color = 255
}
